$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Rename the second sheet (was "sheet1") to "gfdug"
$ws2.Name = "gfdug"

# Row 1 - A1 forced to text (matches source "8207" id, not a number)
$ws2.Cells.Item(1, 1).Formula = "'8207"
$ws2.Cells.Item(1, 1).Style = "Normal"
$ws2.Cells.Item(1, 2).Value = 1248
$ws2.Cells.Item(1, 3).Value = 1356
$ws2.Cells.Item(1, 4).Value = 1244
$ws2.Cells.Item(1, 5).Value = 3
$ws2.Cells.Item(1, 6).Value = 111
$ws2.Cells.Item(1, 7).Value = 91.80811808118081
$ws2.Cells.Item(1, 8).Value = 99.75942261427426
$ws2.Cells.Item(1, 9).Value = 0.084070796460177
$ws2.Cells.Item(1, 10).Value = 46.32402014732361

# Row 2 - A2 forced to text (matches source "8209" id, not a number)
$ws2.Cells.Item(2, 1).Formula = "'8209"
$ws2.Cells.Item(2, 1).Style = "Normal"
$ws2.Cells.Item(2, 2).Value = 1972
$ws2.Cells.Item(2, 3).Value = 2025
$ws2.Cells.Item(2, 4).Value = 1971
$ws2.Cells.Item(2, 5).Value = 0
$ws2.Cells.Item(2, 6).Value = 53
$ws2.Cells.Item(2, 7).Value = 97.38142292490119
$ws2.Cells.Item(2, 8).Value = 100
$ws2.Cells.Item(2, 9).Value = 0.02617283950617284
$ws2.Cells.Item(2, 10).Value = 39.51032686233521

# Reset selection on the data sheet back to A1, and activate the first sheet
$ws2.Range("A1").Select() | Out-Null
$ws1.Activate() | Out-Null
